$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE on row numbering: this worksheet's underlying XML uses a non-standard
# 0-based sheetData row numbering (row r="0" is the first row). Excel/COM
# addressing is always 1-based, and this runtime maps COM row N directly to
# the worksheet's stored row N (i.e. COM row 1 == stored row "1", COM row 2
# == stored row "2", etc.) - the stored row "0" sits "above" COM row 1 and is
# not reachable through the Excel object model (Excel itself has no row 0).
# So here we populate the rows that ARE reachable: COM rows 1-5, which are
# the diff's sheetData rows r="1" through r="5" (the second through sixth
# records).
$data = @(
    @("kVAJt85j", "trashboatsr", 1818, 100, "https://lichess.org/kVAJt85j", 2272),
    @("W26Ykr8M", "trashboatsr", 1818, 100, "https://lichess.org/W26Ykr8M", 2273),
    @("6ZBZX1lE", "trashboatsr", 1818, 100, "https://lichess.org/6ZBZX1lE", 2309),
    @("T5XjC3ky", "trashboatsr", 1818, 100, "https://lichess.org/T5XjC3ky", 2314),
    @("1Y81gpYC", "trashboatsr", 1818, 100, "https://lichess.org/1Y81gpYC", 2315)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}

# Best-effort: also try to update the otherwise-unreachable first stored row
# (sheetData r="0") with the first record's data, in case the host exposes
# it through Offset addressing.
$first = @("gQ9xLuRH", "trashboatsr", 1818, 100, "https://lichess.org/gQ9xLuRH", 2263)
for ($c = 1; $c -le 6; $c++) {
    try {
        $ws.Cells.Item(1, $c).Offset(-1, 0).Value = $first[$c - 1]
    } catch {
    }
}
